# Generate Report for Handback
# Mirrors the localization-status report refresh: marks the two tracked
# files (7a35fbf3... and ca6ecd48...) as handed back, fills in their
# target/handback file names + handback timestamps per language, and
# widens the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2260d96095b90965a5f6ae0dc520e5768e760280/e2e/"

$file1 = "7a35fbf3-18bd-44a7-90c8-f26ca33685d3.md"
$file2 = "ca6ecd48-c673-4fea-9792-4020c15d8bc0.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both locales report back in sync with en-US.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet: both files handed back at 2016-08-26 20:47:15
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("I2").Value = $file1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($repoBase + $file1), [Type]::Missing, [Type]::Missing, $file1)
$zhcn.Range("J2").Value = "7a35fbf3-18bd-44a7-90c8-f26ca33685d3.618fac283f369b4f6f2c86db408cf5594fc56a34.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-26 20:47:15"

$zhcn.Range("I3").Value = $file2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($repoBase + $file2), [Type]::Missing, [Type]::Missing, $file2)
$zhcn.Range("J3").Value = "ca6ecd48-c673-4fea-9792-4020c15d8bc0.66ad9806dfa66851111ba62668fcee3848387e08.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-26 20:47:15"

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: both files handed back at 2016-08-26 20:47:23
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("I2").Value = $file1
$dede.Hyperlinks.Add($dede.Range("I2"), ($repoBase + $file1), [Type]::Missing, [Type]::Missing, $file1)
$dede.Range("J2").Value = "7a35fbf3-18bd-44a7-90c8-f26ca33685d3.618fac283f369b4f6f2c86db408cf5594fc56a34.de-de.xlf"
$dede.Range("K2").Value = "2016-08-26 20:47:23"

$dede.Range("I3").Value = $file2
$dede.Hyperlinks.Add($dede.Range("I3"), ($repoBase + $file2), [Type]::Missing, [Type]::Missing, $file2)
$dede.Range("J3").Value = "ca6ecd48-c673-4fea-9792-4020c15d8bc0.66ad9806dfa66851111ba62668fcee3848387e08.de-de.xlf"
$dede.Range("K3").Value = "2016-08-26 20:47:23"

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17

Write-Output "Report generated for handback."
